$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 8-15: "name" (B) shifts down by two line-entries (line7, line8
# inserted before the extr* rows), and from_bus/to_bus/in_service are recomputed.
$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $false },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A needs the bordered/bold/centered style used throughout column A.
    # Rows 8-15 already carry it; only copy it into the two brand-new rows (16, 17).
    if ($rowNum -gt 15) {
        $ws.Cells.Item($rowNum - 1, 1).Copy()
        $ws.Cells.Item($rowNum, 1).PasteSpecial(-4122)
    }

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
}

$excel.CutCopyMode = $false
